# MOS-23045: Update Master Data as per 22 April Changes
# Add the missing apptyp_code / POA / RNC combinations to the
# master-valid_document sheet and (re)apply the AutoFilter over the
# original data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Apply the AutoFilter to the existing data range *before* the new
#    rows are appended, so the filter (and the hidden _FilterDatabase
#    defined name it implies) stays anchored to A1:G57, matching the
#    state captured right after filtering the original data.
# ---------------------------------------------------------------------
$ws.Range("A1:G57").AutoFilter() | Out-Null

$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "='master-valid_document'!`$A`$1:`$G`$57")
$filterName.Visible = $false

# ---------------------------------------------------------------------
# 2. Give the header row the (number-format) styling that the rest of
#    the sheet already carries: column A keeps the "000" zero-padded
#    numeric format, the remaining header cells line up with the
#    bordered style used throughout the data rows.
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 1).NumberFormat = "000"

# ---------------------------------------------------------------------
# 3. Append the missing POA/RNC (doccat_code/doctyp_code) rows for the
#    apptyp_code values that did not yet have that combination.
# ---------------------------------------------------------------------
$newApptyp = @(3, 4, 7, 8, 11, 12, 15)

$startRow = 58
$i = 0
foreach ($a in $newApptyp) {
    $r = $startRow + $i

    $ws.Cells.Item($r, 1).NumberFormat = "000"
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = "POA"
    $ws.Cells.Item($r, 3).Value = "RNC"
    $ws.Cells.Item($r, 4).Value = "eng"
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"

    $i = $i + 1
}

# ---------------------------------------------------------------------
# 4. Move the active selection to H2, as left by the author after the
#    edit (previously the whole "below the table" block was selected).
# ---------------------------------------------------------------------
$ws.Range("H2").Select() | Out-Null
